{"js": "// Append the trainee's extra \"Wat ben ik tegen gekomen?\" notes to the\n// paragraph that currently ends with \"...daarna de adressen. \", then add a\n// blank paragraph right after it (matches the diff: new runs inserted before\n// the bookmark, plus a following empty <w:p/>).\n\nconst anchorText = \". Het blijkt dat je niet eerst de klant kan verwijderen en daarna de adressen. \";\nconst addition = \"Hoe in database zoeken + string werkt niet omdat mensen dan sqlcode in naamveld kunnen invullen. Hoe paginas alleen toegankelijk te maken voor ingelogde gebruiker (per controller if statements of filters). Kan je ook Hibernate validator een wachtwoord laten controleren? Wachtwoord is niet veilig (via string).\";\n\nconst body = context.document.body;\nconst results = body.search(anchorText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor text not found: \" + anchorText);\n}\n\n// The search hit is the whole sentence; insert the new text right after it,\n// still inside the same paragraph (ahead of the bookmark at the paragraph's\n// tail).\nconst hit = results.items[0];\nhit.insertText(addition, \"After\");\nawait context.sync();\n\n// The paragraph that contains the hit is the one that needs a new, empty\n// paragraph inserted right after it.\nconst paragraph = hit.paragraphs.getFirst();\nparagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n", "ps1": "# Append the trainee's extra \"Wat ben ik tegen gekomen?\" notes to the\n# paragraph that currently ends with \"...daarna de adressen. \", then add a\n# blank paragraph right after it (matches the diff: new text appended before\n# the bookmark, plus a following empty paragraph).\n\n$d = $word.ActiveDocument\n\n$anchorText = \". Het blijkt dat je niet eerst de klant kan verwijderen en daarna de adressen. \"\n$addition = \"Hoe in database zoeken + string werkt niet omdat mensen dan sqlcode in naamveld kunnen invullen. Hoe paginas alleen toegankelijk te maken voor ingelogde gebruiker (per controller if statements of filters). Kan je ook Hibernate validator een wachtwoord laten controleren? Wachtwoord is niet veilig (via string).\"\n\n$rng = $d.Content\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute($anchorText)\n\nif ($found) {\n    # Appending to Range.Text (rather than InsertAfter) keeps the new text\n    # ahead of the paragraph's trailing bookmark, same as the source edit.\n    $rng.Text = $rng.Text + $addition\n\n    # $rng now spans the (extended) sentence; its enclosing paragraph is the\n    # one that needs a new, empty paragraph inserted right after it.\n    $para = $rng.Paragraphs(1)\n    $para.Range.InsertParagraphAfter()\n}\n"}
